$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 48, shifting existing rows 48-54 down to 49-55
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with the "statistic" entry
$ws.Range("A48").Value = "STATO:0000039"
$ws.Range("B48").Value = "statistic"
$ws.Range("C48").Value = "a statistic is a measurement datum to describe a dataset or a variable. It is generated by a calculation on set of observed data."
$ws.Range("D48").Value = "data item"
$ws.Range("E48").Value = ""
$ws.Range("F48").Value = ""
$ws.Range("G48").Value = ""
$ws.Range("H48").Value = ""
$ws.Range("I48").Value = ""
$ws.Range("J48").Value = ""
$ws.Range("K48").Value = ""
$ws.Range("L48").Value = ""
$ws.Range("M48").Value = ""
$ws.Range("N48").Value = ""
$ws.Range("O48").Value = ""
$ws.Range("P48").Value = "LSR 1"
$ws.Range("Q48").Value = "Intervention content and delivery"
$ws.Range("R48").Value = ""
$ws.Range("S48").Value = "External"
$ws.Range("T48").Value = ""
$ws.Range("U48").Value = ""
$ws.Range("V48").Value = "PS"
